# "calculate error of g"
#
# The original sheet had, at rows 9-12:
#   9  average_g   (values only, B9:F9)
#   10 variance    (formulas, B10:F10)
#   11 standard_deviation (values only, B11:F11)
#   12 standard_error     (values only, B12:F12)
#
# The edit inserts 4 new summary rows right after row 9 (g / mean_g /
# standard_deviation / standard_error of g / time), pushing the old
# variance / standard_deviation / standard_error (of distance) rows down
# from 10-12 to 14-16. Row 9 itself keeps its numbers but gains formulas.
#
# Columns J:M (the "distance"/"time" helper table living next to the main
# block) are NOT touched by the insertion - only columns A:F move.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$cols = @("A","B","C","D","E","F")

function Move-Row([string]$fromRow, [string]$toRow) {
    foreach ($col in $cols) {
        $srcCell = $ws.Range($col + $fromRow)
        $dstCell = $ws.Range($col + $toRow)
        if ($srcCell.HasFormula) {
            $dstCell.Formula = $srcCell.Formula
        } else {
            $dstCell.Value = $srcCell.Value()
        }
        $dstCell.HorizontalAlignment = $srcCell.HorizontalAlignment()
        $dstCell.VerticalAlignment = $srcCell.VerticalAlignment()
    }
}

# Relocate old rows 12,11,10 down to 16,15,14 (order doesn't matter - targets
# are currently blank) before we overwrite 10-13 with new content.
Move-Row "12" "16"
Move-Row "11" "15"
Move-Row "10" "14"

# --- Row 9: "average_g" -> "g", add the g = 2*distance/time^2 formulas ---
$ws.Range("A9").Value = "g"
$ws.Range("B9").Formula = "=2*B2/B8^2"
$ws.Range("C9:F9").Formula = "=2*C2/C8^2"

# --- Row 10: mean_g = AVERAGE(g) ---
$ws.Range("A10").Value = "mean_g"
$ws.Range("B10").Formula = "=AVERAGE(B9:F9)"
$ws.Range("B10:F10").Merge()
$ws.Range("B10:F10").HorizontalAlignment = -4108

# --- Row 11: standard_deviation of g = STDEVP(g) ---
$ws.Range("A11").Value = "standard_deviation"
$ws.Range("B11").Formula = "=STDEVP(B9:F9)"
$ws.Range("B11:F11").Merge()
$ws.Range("B11:F11").HorizontalAlignment = -4108

# --- Row 12: standard_error of g = stdev / SQRT(5) ---
$ws.Range("A12").Value = "standard_error"
$ws.Range("B12").Formula = "=B11/SQRT(5)"
$ws.Range("B12:F12").Merge()
$ws.Range("B12:F12").HorizontalAlignment = -4108

# --- Row 13: blank "time" row, centered placeholders ---
$ws.Range("A13").Value = "time"
$ws.Range("B13:F13").HorizontalAlignment = -4108

# Restore column-A label alignment (center/top) on the new rows to match the
# rest of the A column.
foreach ($r in @("9","10","11","12","13")) {
    $ws.Range("A" + $r).HorizontalAlignment = -4108
    $ws.Range("A" + $r).VerticalAlignment = -4160
}

$ws.Range("M11").Select()
